# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" for every cell that
#    holds it (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3).
# 2) Narrow the "zh-cn"/"de-de" status columns (Overview cols E & F, and the
#    "Status" column C on the per-locale sheets) from ~17.22 chars to
#    ~13.41 chars.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: compare with the literal on the left and read via Value2 -
        # PowerShell's -eq coerces the right-hand side to the left-hand
        # side's type, and some status/flag cells hold real booleans, so
        # "$cell.Value2 -eq $oldStatus" would wrongly stringify/compare.
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value2 = $newStatus
        }
    }
}

# Column width: the stored OOXML "width" attribute Excel writes is
# chars + 5/6 (the 5px default cell padding at the workbook's base font),
# rounded to the nearest 1/6 of a character (i.e. nearest pixel). The
# target stored width of 13.4101845877511 sits between the 13.333333...
# and 13.5 pixel-quantized stops; 13.333333... (ColumnWidth ~12.58) is the
# closer of the two, so that's what we set.
$targetColumnWidth = 12.58

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth  # zh-cn
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth  # de-de

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = $targetColumnWidth        # Status

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = $targetColumnWidth        # Status
